# Weekly Fruta/Hortaliza update: insert two new observation rows
# (Granada, Provincia de Curico, week of 45063) above the existing
# data block, pushing the previous rows 12-43 down to 14-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 12, shifting old rows 12:43 -> 14:45
$ws.Rows("12:13").Insert()

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# New row 12
$row12 = @(
    9,
    "Vega Central Mapocho de Santiago",
    "Metropolitana",
    45063,
    13,
    "Fruta",
    100104,
    "Frutos de pepita",
    100104001,
    "Granada",
    "Wonderfull",
    "Especial",
    220,
    10500,
    10500,
    10500,
    "`$/caja 15 kilos granel",
    "Provincia de Curicó",
    700,
    15
)

# New row 13
$row13 = @(
    9,
    "Vega Central Mapocho de Santiago",
    "Metropolitana",
    45063,
    13,
    "Fruta",
    100104,
    "Frutos de pepita",
    100104001,
    "Granada",
    "Wonderfull",
    "Primera",
    250,
    9000,
    9000,
    9000,
    "`$/caja 15 kilos granel",
    "Provincia de Curicó",
    600,
    15
)

for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "12").Value = $row12[$i]
    $ws.Range($cols[$i] + "13").Value = $row13[$i]
}
